# Unit9 text.docx -- insert the word "the" after "MPEG is" so the
# sentence reads "MPEG is the most popular system of the video
# compression. ..." (commit: "added lection6 fort ecology").
#
# Word, when a user types text mid-paragraph, splits the host run at
# the caret and tracks the edit point with the automatic "_GoBack"
# bookmark.  We reproduce that exactly: split off " the" into its own
# run, leave a lone _GoBack bookmark where the caret ended up, and
# drop the _GoBack bookmark that used to sit at the end of the
# paragraph (from the previous edit).

$d = $word.ActiveDocument

# 1) Insert " the" right after "MPEG is".
$insPoint = $d.Content.Duplicate
$insPoint.Find.Execute("MPEG is", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0)
$insPoint.Collapse(0)
$insPoint.InsertAfter(" the")

# 2) Keep " the" as its own run (distinct from the preceding "MPEG is"
#    run) by nudging a character property on and back off -- this
#    forces the engine to keep the run boundary instead of silently
#    re-coalescing it into the neighbouring, identically-formatted run.
$theRange = $d.Content.Duplicate
$theRange.Find.Execute(" the", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0)
$theRange.Font.Bold = 1
$theRange.Font.Bold = 0

# 3) Likewise keep the remaining tail of the original run (" most
#    popular system ... P-frame store") separate from the run that
#    follows it in the document ("s only the difference ...").
$tailText = " most popular system of the video compression. This system splits the single data stream into video and audio with the different algorithms. There are some types of frames for this method. First one is Intra Frame (or I-frame), which is compressed using the picture itself like JPEG. The second one is predicted frame (or P-frame). P-frame store"
$tailRange = $d.Content.Duplicate
$tailRange.Find.Execute($tailText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$tailRange.Font.Bold = 1
$tailRange.Font.Bold = 0

# 4) Drop the _GoBack bookmark that used to be at the end of the
#    paragraph (leftover from whatever edit happened before this one).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 5) Re-create _GoBack right after the newly typed " the", matching
#    where Word would leave it after this edit.
$goBackPoint = $d.Content.Duplicate
$goBackPoint.Find.Execute(" the", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
$goBackPoint.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBackPoint)
